# Sync attendance_reports: normalize the "Recorded By" (column G) author
# list ordering. Each cell holds a comma-separated list of recorder
# identifiers (e.g. "System, dnasr281@gmail.com"); re-sort that list
# case-insensitively so "System" no longer always leads, while leaving
# single-value cells and cells referencing admin@admin.com untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    # Only touch multi-value lists; leave single entries alone.
    if (-not ($val -like "*, *")) {
        continue
    }

    # Rows that mention admin@admin.com were left untouched in the sync.
    if ($val.Contains("admin@admin.com")) {
        continue
    }

    $parts = $val -split ", "
    $sorted = $parts | Sort-Object { $_.ToLower() }
    $newVal = $sorted -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
